$p = $ppt.ActivePresentation

# --- Slide 1: "Lectures" -> "Information meeting" (bold green run) ---
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange
$para1 = $tr1.Paragraphs(10, 1)
$ptext1 = $para1.Text
$pos1 = $ptext1.IndexOf("Lectures")
$runStart1 = $para1.Start + $pos1
$runChars1 = $tr1.Characters($runStart1, 8)
$runChars1.Text = "Information meeting"

# --- Slide 12: merge "Vigtige " + "meddelelser" runs into one run ---
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(2)
$tr12 = $sh12.TextFrame.TextRange
$para2 = $tr12.Paragraphs(2, 1)
$ptext2 = $para2.Text
$pos2 = $ptext2.IndexOf("Vigtige")
$runStart2 = $para2.Start + $pos2
$runChars2 = $tr12.Characters($runStart2, 19)
$runChars2.Text = "Vigtige meddelelser"

# --- Slide 12: merge "Postings " + "on the "Discussion forum"" runs into one run ---
$para3 = $tr12.Paragraphs(3, 1)
$t3 = $para3.Text
$trimmed3 = $t3.TrimEnd([char]13)
$chars3 = $tr12.Characters($para3.Start, $trimmed3.Length)
$chars3.Text = "Postings on the “Discussion forum”"
